$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update a few existing daily totals in May/2025 (rows 7, 18, 19)
$ws.Range("B7").Value = 38247.95
$ws.Range("B18").Value = 25760.86
$ws.Range("B19").Value = 11480.51

# Insert a new daily record (Dia 28) for May/2025 right after the current
# last May row (row 19), shifting all subsequent rows down by one.
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = 28
$ws.Range("B20").Value = 1534.16
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 2025
$ws.Range("E20").Value = "05/2025"
